# Applies the edits described in the diff to the A35 sheet of Statistics.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("A35")

# Row 9: K9 changes from 4 to 5
$ws.Range("K9").Value = 5

# Row 9: O9 loses its shared formula and becomes a static value of 14
$ws.Range("O9").Value = 14

# Row 9: P9 keeps its formula (=O9 via shared formula) -> recalculates to 14 automatically

# Row 11: L11 changes from 192 to 193
$ws.Range("L11").Value = 193

# Row 15: L15 keeps formula SUM(L3:L14) -> recalculates to 1476 automatically

# Update the selected cell/sqref shown in the saved sheet view
$ws.Range("D17").Select()

$excel.CalculateFullRebuild()
